# Add the new weekly ranking sheet "2026-01-28" at the end of the workbook,
# populate it with rank/title/volume/publisher data, and re-apply the
# yellow "low volume count" highlight (fillId=2 in styles.xml) to the same
# cells that carry it in the source data — done by copying the format from
# an existing styled cell elsewhere in the workbook, so no duplicate fill
# gets registered in the style table.

$wb = $excel.ActiveWorkbook

$sheetName = "2026-01-28"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = $sheetName

# Pipe-delimited rank|title|volume|publisher rows (row 1 is the header).
# Literal here-string: avoids PowerShell variable expansion on "$" etc.
$rowText = @'
rank|title|volume|publisher
1|ブルーロック|37|
2|シャングリラ・フロンティア ~クソゲーハンター、神ゲーに挑まんとす~|25|
3|はじめの一歩|145|
4|金色のガッシュ!! 2 Page|36|
5|彼女、お借りします|44|
6|フェルマーの料理|7|
7|金色のガッシュ!!|6|
8|カッコウの許嫁|31|
9|葬送のフリーレン|15|
10|黒岩メダカに私の可愛いが通じない|22|
11|ダンダダン|22|
12|TSUYOSHI 誰も勝てない、アイツには|29|
13|うるわしの宵の月|10|
14|ビジネス婚ー好きになったら離婚しますー|1|
15|僕の心のヤバイやつ|13|
16|婚約者が浮気しているようなんですけど私は流行りの悪役令嬢ってことであってますか?|3|
17|波うららかに、めおと日和|10|
18|ONE PIECE|113|
19|#神奈川に住んでるエルフ|1|
20|追放された没落令嬢は拳ひとつで異世界を生き延びる! コミック版|1|
21|すだちの魔王城|12|
22|弱虫ペダル|98|
23|スキル? ねぇよそんなもん! ~不遇者たちの才能開花~ コミック版|1|
24|異常死体解剖ファイル|1|
25|契約恋人は、前世で私を裏切った男です1|1|
26|おしかけ勇者嫁 勇者は放逐されたおっさんを追いかけ、スローライフを応援する コミック版|1|
27|異常死体解剖ファイル|2|
28|ケジメつけさせてもらいます。元ヤン弁護士 東矢斎|6|
29|静かなるドン ― もうひとつの最終章 ―|7|
30|転生貴族、鑑定スキルで成り上がる ~弱小領地を受け継いだので、優秀な人材を増やしていたら、最強領地になってた~|20|
31|WORST外伝 ゼットン先生|17|
32|ヘルモード ~やり込み好きのゲーマーは廃設定の異世界で無双する~はじまりの召喚士13|13|
33|恋せよまやかし天使ども|6|
34|龍と苺|24|
35|カグラバチ|10|
36|呪術廻戦≡(モジュロ)|1|
37|わたし、今日から「おひとりさま」|4|
38|国民的アイドルが弟になったら|10|
39|ジャンク・ランク・ファミリー|20|
40|転生幼女はあきらめない|11|
41|貴族転生 ~恵まれた生まれから最強の力を得る~|10|
42|旅はオールパス|1|
43|イジメ島第1話|1|
44|リベンジノート|6|
45|ケジメつけさせてもらいます。元ヤン弁護士 東矢斎|4|
46|バニシング・ツイン~私の中の君~|3|
47|リセット~もしも夫以外の誰かと~|3|
48|義兄の純愛~初めての恋もカラダも、エリート弁護士に教えられました~|1|
49|魔法歌姫マジカルギンガ 第27話|27|
50|ザ・ファブル The third secret|3|
51|銭(インチキ)の力で、戦国の世を駆け抜ける。|8|
52|片田舎のおっさん、剣聖になる~ただの田舎の剣術師範だったのに、大成した弟子たちが俺を放ってくれない件~|8|
53|世界に一人、全属性魔法の使い手|8|
54|穏やか貴族の休暇のすすめ。@COMIC|14|
55|七つ屋志のぶの宝石匣|26|
56|なつめとなつめ|1|
57|リベンジノート|4|
58|リベンジノート|5|
59|ケジメつけさせてもらいます。元ヤン弁護士 東矢斎|5|
60|国民的アイドルが弟になったら|6|
61|国民的アイドルが弟になったら|7|
62|国民的アイドルが弟になったら|8|
63|国民的アイドルが弟になったら|9|
64|次に買うマンガ、この1話で決めよう! ~王道こそ至高にして最強! 王道恋愛編~|4|
65|四天王最弱の自立計画(コミック) 1話|1|
66|氷の侯爵令嬢は、魔狼騎士に甘やかに溶かされる|1|
67|氷の侯爵令嬢は、魔狼騎士に甘やかに溶かされる|2|
68|氷の侯爵令嬢は、魔狼騎士に甘やかに溶かされる|3|
69|四獣封地伝|1|
70|期限つき皇女のはずが、うまくやりすぎてしまったようです|1|
71|期限つき皇女のはずが、うまくやりすぎてしまったようです|2|
72|期限つき皇女のはずが、うまくやりすぎてしまったようです|3|
73|聖女と公爵様の晩酌 ~前世グルメで餌付けして、のんびり楽しい偽物夫婦ぐらし~|1|
74|にぶんのいち夫婦|1|
75|俺だけレベルアップな件|23|
76|お気楽領主の楽しい領地防衛|1|
77|お気楽領主の楽しい領地防衛|4|
78|お気楽領主の楽しい領地防衛|7|
79|追放された転生重騎士はゲーム知識で無双する|16|
80|異世界でテイムした最強の使い魔は、幼馴染の美少女でした|6|
81|異世界じゃスローライフはままならない ~聖獣の主人は島育ち~6|6|
82|お気楽領主の楽しい領地防衛|3|
83|お気楽領主の楽しい領地防衛|5|
84|お気楽領主の楽しい領地防衛|6|
85|るろうに剣心―明治剣客浪漫譚・北海道編―|10|
86|葬送のフリーレン|14|
87|呪術廻戦|30|
88|呪術廻戦|29|
89|お気楽領主の楽しい領地防衛|2|
90|旅はオールパス|2|
91|旅はオールパス|3|
92|40歳捨てられ花嫁、なぜか年下副社長に溺愛されてます1|1|
93|勇者の出番ねぇからっ!! ~異世界転生するけど俺は脇役と言われました~ コミック版|1|
94|俺はこの世界がモブでもになれば最強になれることを知っている@COMIC 第1話|1|
95|乙女ゲームヒロインの『引き立て役の妹』に転生したので立場を奪ってやることにした。@COMIC 第1話|1|
96|死の運命を回避するために、未来の大公様、私と結婚してください! 第1話|1|
97|売られた聖女は異郷の王の愛を得る 第1話|1|
98|アマーリエと悪食公爵~孤独な令嬢は心のすべてを食べられたい~ 1皿|1|
99|起きたら20年後なんですけど! ~悪役令嬢のその後のその後~ 1(アリアンローズコミックス)|1|
100|「《邪神の血》が流れている」と言われ、神聖教会を追放された神父です。 ~理不尽な理由で教会を追い出されたら、信仰対象の女神様も一緒についてきちゃいました~ コミック版|1|
'@

$lines = $rowText -split "`n"
$rowCount = $lines.Count
$colCount = 4

$data = New-Object 'object[,]' $rowCount,$colCount

for ($r = 0; $r -lt $rowCount; $r++) {
    $fields = $lines[$r] -split "\|"
    $rank = $fields[0]
    $title = $fields[1]
    $volume = $fields[2]
    $publisher = $fields[3]

    if ($r -eq 0) {
        # header row: keep every column as text, matching the source sheet
        $data[$r,0] = $rank
        $data[$r,1] = $title
        $data[$r,2] = $volume
        $data[$r,3] = $publisher
    } else {
        $data[$r,0] = [double]$rank
        $data[$r,1] = $title
        $data[$r,2] = [double]$volume
        $data[$r,3] = $null
    }
}

$endCell = $newSheet.Cells.Item($rowCount, $colCount)
$targetRange = $newSheet.Range($newSheet.Cells.Item(1,1), $endCell)
$targetRange.Value = $data

# Re-apply the yellow "low volume number" highlight style to the same
# volume cells that had it in the previous weekly sheets, by copying the
# cell format from an existing highlighted cell (reuses the existing fill
# instead of registering a new one).
$styleSource = $wb.Worksheets.Item("2026-01-21").Range("C10")
$styleSource.Copy()

$highlightRows = @(15,17,20,21,24,25,26,27,28,37,43,44,47,48,49,51,57,66,67,68,69,70,71,72,73,74,75,77,83,90,91,92,93,94,95,96,97,98,99,100,101)
foreach ($r in $highlightRows) {
    $newSheet.Range("C" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
